$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule-P1")

$ws.Range("G6").Value = 42517
$ws.Range("G7").Value = 42517
$ws.Range("G8").Value = 42517
$ws.Range("G9").Value = 42517
